$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 89 (date + high corrected by the R script) ---
$ws.Range("A89").Value = 45470.2916666667
$ws.Range("C89").Value = 3.22000002861023

# --- Append new row 90 with the latest scraped OHLCV data ---

# Column A (date) carries the custom "yyyy-mm-dd hh:mm:ss" style used by the
# rest of the column; copy that formatting from the row above first, then set
# the value so the cell lands on the same style (s="1") instead of minting a
# new one.
$ws.Range("A89").Copy() | Out-Null
$ws.Range("A90").PasteSpecial(-4122) | Out-Null
$ws.Range("A90").Value = 45471.6493865741

# Plain numeric columns: B (volume), C (high), D (low), E (open), F (close)
$ws.Range("B90").Value = 198000
$ws.Range("C90").Value = 3.45000004768372
$ws.Range("D90").Value = 3.0699999332428
$ws.Range("E90").Value = 3.33999991416931
$ws.Range("F90").Value = 3.29999995231628

# Columns G (adj_close) and H (ticker) store their values as shared strings
# (text), not numbers, matching every other row in the sheet. Force text
# entry via a temporary "@" number format, then restore the plain
# (General/unstyled) formatting used elsewhere in those columns by copying it
# from row 2.
$ws.Range("G90:H90").NumberFormat = "@"
$ws.Range("G90").Value = "3.29999995231628"
$ws.Range("H90").Value = "ESPE.MI"

$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("G90:H90").PasteSpecial(-4122) | Out-Null
